$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns D (Price) and E (Volume) are stored as text in the source sheet.
# For values that look like plain numbers (e.g. "244.07", "1.000"), prefix
# with a literal apostrophe -- the same quote-prefix Excel itself uses to
# force text entry -- so they are not silently re-typed as numbers and lose
# formatting (trailing zeros, trailing-dot-grouping, scientific notation, etc).

$ws.Range('D2').Value = '30.355.95'
$ws.Range('E2').Value = '  -0.66%  '
$ws.Range('D3').Value = '1.871.72'
$ws.Range('E3').Value = '  -0.59%  '
$ws.Range('E4').Value = '  -0.01%  '
$ws.Range('D5').Value = "'244.07"
$ws.Range('E5').Value = '  -2.16%  '
$ws.Range('E6').Value = '  -0.02%  '
$ws.Range('D7').Value = "'0.4717"
$ws.Range('E7').Value = '  -0.63%  '
$ws.Range('D8').Value = "'0.2877"
$ws.Range('E8').Value = '  -1.92%  '
$ws.Range('D9').Value = "'0.06483"
$ws.Range('E9').Value = '  -0.72%  '
$ws.Range('D10').Value = "'21.73"
$ws.Range('E10').Value = '  -1.34%  '
$ws.Range('D11').Value = "'0.07798"
$ws.Range('E11').Value = '  +0.86%  '
$ws.Range('D12').Value = '1.872.36'
$ws.Range('E12').Value = '  -0.55%  '
$ws.Range('D13').Value = "'96.12"
$ws.Range('D14').Value = "'0.7152"
$ws.Range('E14').Value = '  -4.20%  '
$ws.Range('D15').Value = "'5.133"
$ws.Range('E15').Value = '  -2.78%  '
$ws.Range('D16').Value = "'283.14"
$ws.Range('E16').Value = '  +2.55%  '
$ws.Range('D17').Value = '30.351.45'
$ws.Range('E17').Value = '  -1.17%  '
$ws.Range('D18').Value = "'12.97"
$ws.Range('E18').Value = '  -2.02%  '
$ws.Range('D19').Value = "'1.000"
$ws.Range('E19').Value = '  -0.02%  '
$ws.Range('D20').Value = "'0.000007501"
$ws.Range('E20').Value = '  -0.57%  '
$ws.Range('D21').Value = '2.114.35'
$ws.Range('E21').Value = '  -0.63%  '
$ws.Range('B22').Value = 'Uniswap'
$ws.Range('C22').Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range('D22').Value = "'5.287"
$ws.Range('E22').Value = '  -0.52%  '
$ws.Range('B23').Value = 'BinanceUSD'
$ws.Range('C23').Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range('D23').Value = "'1.000"
$ws.Range('E23').Value = '  -0.04%  '
$ws.Range('D24').Value = "'6.257"
$ws.Range('E24').Value = '  +0.24%  '
$ws.Range('D25').Value = "'163.15"
$ws.Range('E25').Value = '  -0.78%  '
$ws.Range('D26').Value = "'9.020"
$ws.Range('E26').Value = '  -2.06%  '
$ws.Range('D27').Value = "'18.74"
$ws.Range('E27').Value = '  -0.74%  '
$ws.Range('D28').Value = "'1.884"
$ws.Range('E28').Value = '  -1.62%  '
$ws.Range('D29').Value = "'1.317"
$ws.Range('E29').Value = '  -1.97%  '
$ws.Range('D30').Value = "'0.09595"
$ws.Range('E30').Value = '  -1.03%  '
$ws.Range('D31').Value = "'1.488"
$ws.Range('E31').Value = '  -1.10%  '
$ws.Range('D32').Value = "'4.222"
$ws.Range('E32').Value = '  -1.83%  '
$ws.Range('D33').Value = "'4.126"
$ws.Range('E33').Value = '  -0.78%  '
$ws.Range('D34').Value = "'0.04836"
$ws.Range('E34').Value = '  -0.53%  '
$ws.Range('E35').Value = '  -0.51%  '
$ws.Range('D36').Value = "'0.6892"
$ws.Range('E36').Value = '  -1.30%  '
$ws.Range('D37').Value = "'2.710"
$ws.Range('E37').Value = '  -0.36%  '
$ws.Range('D38').Value = "'0.01892"
$ws.Range('E38').Value = '  -0.24%  '
$ws.Range('D39').Value = "'2.821"
$ws.Range('E39').Value = '  +1.75%  '
$ws.Range('D40').Value = "'75.47"
$ws.Range('E40').Value = '  +1.56%  '
$ws.Range('D41').Value = "'6.231"
$ws.Range('E41').Value = '  -1.18%  '
$ws.Range('B42').Value = 'RenderToken'
$ws.Range('C42').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D42').Value = "'1.921"
$ws.Range('E42').Value = '  -4.82%  '
$ws.Range('B43').Value = 'TheSandbox'
$ws.Range('C43').Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range('D43').Value = "'0.4195"
$ws.Range('E43').Value = '  -1.07%  '
$ws.Range('E44').Value = '  -0.07%  '
$ws.Range('D45').Value = "'0.8239"
$ws.Range('E45').Value = '  -1.96%  '
$ws.Range('D46').Value = "'100.77"
$ws.Range('E46').Value = '  -1.34%  '
$ws.Range('D47').Value = "'9.773"
$ws.Range('E47').Value = '  +3.00%  '
$ws.Range('D48').Value = "'7.004"
$ws.Range('E48').Value = '  -0.50%  '
$ws.Range('D49').Value = "'35.17"
$ws.Range('E49').Value = '  -1.44%  '
$ws.Range('D50').Value = "'899.34"
$ws.Range('E50').Value = '  -1.57%  '
$ws.Range('D51').Value = "'0.05735"
$ws.Range('E51').Value = '  +0.18%  '
